# Mass Interview uncheck for Only Available Interviewers
# Adds new LIVE interview history rows to the AMSIN sheet (rows 22 & 23,
# plus a style fix on row 21) and one new row to the BETA sheet (row 11).

$wb = $excel.ActiveWorkbook
$wsAmsin = $wb.Worksheets.Item("AMSIN")
$wsBeta  = $wb.Worksheets.Item("BETA")

# ---------------------------------------------------------------------
# 1) BETA!A11/C11/E11/F11/G11 need to end up with NO explicit style
#    (default style) in the final file. The only cells in the whole
#    workbook that currently carry that default/unstyled formatting are
#    AMSIN!A21/C21/D21/E21/F21/G21 - grab a copy of that formatting
#    *before* we touch row 21 below.
# ---------------------------------------------------------------------
$wsAmsin.Range("A21").Copy()
$wsBeta.Range("A11").PasteSpecial(-4122)   # xlPasteFormats

$wsAmsin.Range("C21").Copy()
$wsBeta.Range("C11").PasteSpecial(-4122)

$wsAmsin.Range("D21").Copy()
$wsBeta.Range("E11").PasteSpecial(-4122)

$wsAmsin.Range("D21").Copy()
$wsBeta.Range("F11").PasteSpecial(-4122)

$wsAmsin.Range("D21").Copy()
$wsBeta.Range("G11").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 2) BETA row 11 values (A/C = text, forced via Text format then
#    reverted to keep the pasted "no style" formatting; E/F/G numeric;
#    B becomes a date/time value with style 10; D keeps its existing
#    style 9 and simply gains a value)
# ---------------------------------------------------------------------
$wsBeta.Range("A11").NumberFormat = "@"
$wsBeta.Range("A11").Value = "2021-10-28"
$wsAmsin.Range("A21").Copy()
$wsBeta.Range("A11").PasteSpecial(-4122)

$wsBeta.Range("C11").NumberFormat = "@"
$wsBeta.Range("C11").Value = "152_beta"
$wsAmsin.Range("C21").Copy()
$wsBeta.Range("C11").PasteSpecial(-4122)

$wsBeta.Range("B11").Value = 44497.62992095913
$wsBeta.Range("B11").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$wsBeta.Range("D11").Value = 105
$wsBeta.Range("E11").Value = 105
$wsBeta.Range("F11").Value = 0
$wsBeta.Range("G11").Value = 3.31

# ---------------------------------------------------------------------
# 3) AMSIN row 21: apply the "data row" style (copied from row 20) to
#    A/C/D/E/F/G, and refresh B21's timestamp value.
# ---------------------------------------------------------------------
$wsAmsin.Range("A20").Copy()
$wsAmsin.Range("A21").PasteSpecial(-4122)

$wsAmsin.Range("C20").Copy()
$wsAmsin.Range("C21").PasteSpecial(-4122)

$wsAmsin.Range("D20").Copy()
$wsAmsin.Range("D21").PasteSpecial(-4122)

$wsAmsin.Range("E20").Copy()
$wsAmsin.Range("E21").PasteSpecial(-4122)

$wsAmsin.Range("F20").Copy()
$wsAmsin.Range("F21").PasteSpecial(-4122)

$wsAmsin.Range("G20").Copy()
$wsAmsin.Range("G21").PasteSpecial(-4122)

$wsAmsin.Range("B21").Value = 44475.70902158565

# ---------------------------------------------------------------------
# 4) AMSIN row 22 (new): full data row, same style as row 20.
# ---------------------------------------------------------------------
$wsAmsin.Range("A20").Copy()
$wsAmsin.Range("A22").PasteSpecial(-4122)
$wsAmsin.Range("A22").NumberFormat = "@"
$wsAmsin.Range("A22").Value = "2021-10-26"
$wsAmsin.Range("A20").Copy()
$wsAmsin.Range("A22").PasteSpecial(-4122)

$wsAmsin.Range("C20").Copy()
$wsAmsin.Range("C22").PasteSpecial(-4122)
$wsAmsin.Range("C22").NumberFormat = "@"
$wsAmsin.Range("C22").Value = "152_fstcycle"
$wsAmsin.Range("C20").Copy()
$wsAmsin.Range("C22").PasteSpecial(-4122)

$wsAmsin.Range("D20").Copy()
$wsAmsin.Range("D22").PasteSpecial(-4122)
$wsAmsin.Range("D22").Value = 105

$wsAmsin.Range("E20").Copy()
$wsAmsin.Range("E22").PasteSpecial(-4122)
$wsAmsin.Range("E22").Value = 104

$wsAmsin.Range("F20").Copy()
$wsAmsin.Range("F22").PasteSpecial(-4122)
$wsAmsin.Range("F22").Value = 1

$wsAmsin.Range("G20").Copy()
$wsAmsin.Range("G22").PasteSpecial(-4122)
$wsAmsin.Range("G22").Value = 3.63

$wsAmsin.Range("B22").Value = 44495.65370459491
$wsAmsin.Range("B22").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# ---------------------------------------------------------------------
# 5) AMSIN row 23 (new): full data row, same style as row 20.
# ---------------------------------------------------------------------
$wsAmsin.Range("A20").Copy()
$wsAmsin.Range("A23").PasteSpecial(-4122)
$wsAmsin.Range("A23").NumberFormat = "@"
$wsAmsin.Range("A23").Value = "2021-10-27"
$wsAmsin.Range("A20").Copy()
$wsAmsin.Range("A23").PasteSpecial(-4122)

$wsAmsin.Range("C20").Copy()
$wsAmsin.Range("C23").PasteSpecial(-4122)
$wsAmsin.Range("C23").NumberFormat = "@"
$wsAmsin.Range("C23").Value = "152_scndcycle"
$wsAmsin.Range("C20").Copy()
$wsAmsin.Range("C23").PasteSpecial(-4122)

$wsAmsin.Range("D20").Copy()
$wsAmsin.Range("D23").PasteSpecial(-4122)
$wsAmsin.Range("D23").Value = 105

$wsAmsin.Range("E20").Copy()
$wsAmsin.Range("E23").PasteSpecial(-4122)
$wsAmsin.Range("E23").Value = 93

$wsAmsin.Range("F20").Copy()
$wsAmsin.Range("F23").PasteSpecial(-4122)
$wsAmsin.Range("F23").Value = 12

$wsAmsin.Range("G20").Copy()
$wsAmsin.Range("G23").PasteSpecial(-4122)
$wsAmsin.Range("G23").Value = 6.95

$wsAmsin.Range("B23").Value = 44496.66827148148
$wsAmsin.Range("B23").NumberFormat = "YYYY-MM-DD HH:MM:SS"
